$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Line")

# --- Update H6 and H7 values ---
$ws.Range("H6").Value = 2900
$ws.Range("H7").Value = -2900

# --- Add new shared string "EE00-FI00" via row 20/21 data ---
$ws.Range("A20").Value = "EE00-FI00"
$ws.Range("B20").Value = "Export Capacity"
$ws.Range("C20").Value = "Interconnection"
$ws.Range("D20").Value = "Distributed Energy"
$ws.Range("E20").Value = "Reference Grid"
$ws.Range("F20").Value = 2040
$ws.Range("G20").Value = 1984
$ws.Range("H20").Value = 1700

$ws.Range("A21").Value = "EE00-FI00"
$ws.Range("B21").Value = "Import Capacity"
$ws.Range("C21").Value = "Interconnection"
$ws.Range("D21").Value = "Distributed Energy"
$ws.Range("E21").Value = "Reference Grid"
$ws.Range("F21").Value = 2040
$ws.Range("G21").Value = 1984
$ws.Range("H21").Value = -1700

# --- Update selection (active cell) ---
$ws.Range("I8").Select()
